# Generate Report for Handback
# The handback for "9c788acf-9ee0-447c-a4ae-d1a226a7d50f.md" has completed:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The per-locale "Latest Handback DateTime" is refreshed to the handback time
#   - The per-locale "Error Detail" (stale-handback warning) is cleared
# This mirrors the same update across the Overview sheet and each locale sheet
# (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "2016-08-19 06:46:57"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "2016-08-19 06:47:12"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
